$wb = $excel.ActiveWorkbook

# --- Sheet 2: "anel de retenção" ---------------------------------------
$ws2 = $wb.Worksheets.Item("anel de retenção")

# Insert a blank row above row 12 -> pushes the data block (rows 12-33)
# down to (rows 13-34).
$ws2.Rows("12:12").Insert()

# Move/resize the ring picture (same size, new position).
$pic = $ws2.Shapes.Item("Figura 1")
$pic.Left = 3226320 / 12700
$pic.Top = 91080 / 12700
$pic.Width = 2547000 / 12700
$pic.Height = 1863360 / 12700

# Selection / view state for sheet 2
$ws2.Range("K8").Select()
$excel.ActiveWindow.ScrollRow = 16

# --- Sheet 1: "data_csv" -------------------------------------------------
$ws1 = $wb.Worksheets.Item("data_csv")
$ws1.Range("S32").Select()
